$d = $word.ActiveDocument

$replacements = @(
    @("22×39=", "79×88="),
    @("28×57=", "13×90="),
    @("89×86=", "66×84="),
    @("81×63=", "76×70="),
    @("15×89=", "95×89="),
    @("34×65=", "95×15="),
    @("84×55=", "56×20="),
    @("51×15=", "67×46="),
    @("21×89=", "36×80="),
    @("30×44=", "89×18="),
    @("60×33=", "40×41="),
    @("48×54=", "53×47="),
    @("31×99=", "46×88="),
    @("49×61=", "91×46="),
    @("28×53=", "52×17="),
    @("85×58=", "80×85="),
    @("72×94=", "26×83="),
    @("76×82=", "39×86="),
    @("62×43=", "19×39="),
    @("90×91=", "81×73="),
    @("97×66=", "27×19="),
    @("60×26=", "36×38="),
    @("95×48=", "82×49="),
    @("84×78=", "49×22="),
    @("59×74=", "55×74=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
